# Fruta / hortaliza, semanal
# Inserts a new weekly price record for "Acelga" (Agrícola del Norte S.A. de
# Arica) as row 50 of the data sheet, pushing the existing rows 50-71 down
# to rows 51-72 (entire used range grows from A1:R71 to A1:R72).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 50 downward (through the end of the sheet) by one row, so the
# existing row 50 data ends up in row 51, etc.
$ws.Rows.Item(50).Insert()

# Populate the newly inserted row 50 with the new weekly record.
$ws.Cells.Item(50, 1).Value = 1
$ws.Cells.Item(50, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(50, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(50, 4).Value = 44806
$ws.Cells.Item(50, 5).Value = 15
$ws.Cells.Item(50, 6).Value = 100112009
$ws.Cells.Item(50, 7).Value = "Acelga"
$ws.Cells.Item(50, 8).Value = "Sin especificar"
$ws.Cells.Item(50, 9).Value = "Primera"
$ws.Cells.Item(50, 10).Value = 200
$ws.Cells.Item(50, 11).Value = 1300
$ws.Cells.Item(50, 12).Value = 1500
$ws.Cells.Item(50, 13).Value = 1400
$ws.Cells.Item(50, 14).Value = "$/atado 2,5 a 3 kilos"
$ws.Cells.Item(50, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(50, 16).Value = 467
$ws.Cells.Item(50, 17).Value = 3
$ws.Cells.Item(50, 18).Value = "Hortaliza"
